$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "61÷3=20, 1"
$tbl.Cell(1, 2).Range.Text = "39÷7=5, 4"
$tbl.Cell(1, 3).Range.Text = "32÷3=10, 2"
$tbl.Cell(1, 4).Range.Text = "55÷8=6, 7"
$tbl.Cell(1, 5).Range.Text = "73÷2=36, 1"

$tbl.Cell(5, 1).Range.Text = "29÷3=9, 2"
$tbl.Cell(5, 2).Range.Text = "40÷3=13, 1"
$tbl.Cell(5, 3).Range.Text = "96÷9=10, 6"
$tbl.Cell(5, 4).Range.Text = "54÷6=9, 0"
$tbl.Cell(5, 5).Range.Text = "37÷8=4, 5"

$tbl.Cell(9, 1).Range.Text = "93÷8=11, 5"
$tbl.Cell(9, 2).Range.Text = "76÷2=38, 0"
$tbl.Cell(9, 3).Range.Text = "28÷2=14, 0"
$tbl.Cell(9, 4).Range.Text = "94÷6=15, 4"
$tbl.Cell(9, 5).Range.Text = "86÷4=21, 2"

$tbl.Cell(13, 1).Range.Text = "15÷7=2, 1"
$tbl.Cell(13, 2).Range.Text = "38÷8=4, 6"
$tbl.Cell(13, 3).Range.Text = "43÷7=6, 1"
$tbl.Cell(13, 4).Range.Text = "61÷6=10, 1"
$tbl.Cell(13, 5).Range.Text = "29÷2=14, 1"

$tbl.Cell(17, 1).Range.Text = "32÷4=8, 0"
$tbl.Cell(17, 2).Range.Text = "67÷2=33, 1"
$tbl.Cell(17, 3).Range.Text = "26÷8=3, 2"
$tbl.Cell(17, 4).Range.Text = "59÷3=19, 2"
$tbl.Cell(17, 5).Range.Text = "78÷8=9, 6"

Write-Output "done"